$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("C7:H8").ClearContents()
    $ws.Range("C15:H16").ClearContents()
}
